# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-09-03 (serial 45172) to 2023-09-06 (serial 45175).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$newDate = 45175
$firstRow = 2
$lastRow = 295

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
